$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the standalone "7455355 - Robson da Silva Rocha" value row
# that sits under "Docentes responsáveis:" with no label in column A).
# Everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# Row 10 ("Objetivos:") B/C now hold the responsible professor info instead
# of the old objectives paragraph.
$ws.Range("B10:C10").Value = "7455355 - Robson da Silva Rocha"

# Row 13 ("Programa resumido:", after the shift) B/C now just say "Semestral".
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 ("Programa:", after the shift) B/C now hold the activation date.
$ws.Range("B15:C15").Value = "01/01/2022"

# Row 18 ("Método:", after the shift) B/C now hold the responsible professor
# info again.
$ws.Range("B18:C18").Value = "7455355 - Robson da Silva Rocha"
